$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 295; all existing rows 295-350 shift down to 296-351.
$ws.Rows.Item(295).Insert()

# Populate the newly inserted row 295 with the new weekly price record.
$ws.Cells.Item(295, 1).Value = 9
$ws.Cells.Item(295, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(295, 3).Value = "Metropolitana"
$ws.Cells.Item(295, 4).Value = 45005
$ws.Cells.Item(295, 5).Value = 13
$ws.Cells.Item(295, 6).Value = 100112001
$ws.Cells.Item(295, 7).Value = "Berenjena"
$ws.Cells.Item(295, 8).Value = "Sin especificar"
$ws.Cells.Item(295, 9).Value = "Primera"
$ws.Cells.Item(295, 10).Value = 90
$ws.Cells.Item(295, 11).Value = 7000
$ws.Cells.Item(295, 12).Value = 8000
$ws.Cells.Item(295, 13).Value = 7500
$ws.Cells.Item(295, 14).Value = "$/caja 50 unidades"
$ws.Cells.Item(295, 15).Value = "Región Metropolitana"
$ws.Cells.Item(295, 16).Value = 150
$ws.Cells.Item(295, 17).Value = 50
$ws.Cells.Item(295, 18).Value = "Hortaliza"
